# Add new German/Spanish/French (and a few Portuguese) translations to the
# "Worksheet" and "Sheet1" sheets of the multilingual workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Worksheet": columns are
#   A = identifiers, B = English, C = German, D = Portuguese,
#   E = Spanish, F = French
# Row 1 = headers, Row 2 ("Welcome!") is already fully translated.
# Rows 3-18 get new German / Spanish / French text (Portuguese stays a
# copy of English, as it already was).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Worksheet")

# Row 3 - Save as...
$ws1.Range("C3").Value = "Speichern als..."
$ws1.Range("E3").Value = "olianore"
$ws1.Range("F3").Value = "Enregistrer sous..."

# Row 4 - View
$ws1.Range("C4").Value = "Ansicht"
$ws1.Range("E4").Value = "Ver"
$ws1.Range("F4").Value = "cava"

# Row 5 - About...
$ws1.Range("C5").Value = "Über..."
$ws1.Range("E5").Value = "Acerca de..."
$ws1.Range("F5").Value = "À propos..."

# Row 6 - About Me
$ws1.Range("C6").Value = "Über mich"
$ws1.Range("E6").Value = "Acerca de mí"
$ws1.Range("F6").Value = "À propos de moi"

# Row 7 - Current Password
$ws1.Range("C7").Value = "Aktuelles Passwort"
$ws1.Range("E7").Value = "Contraseña actual"
$ws1.Range("F7").Value = "Mot de passe actuel"

# Row 8 - New Password
$ws1.Range("C8").Value = "Neues Passwort"
$ws1.Range("E8").Value = "Nueva contraseña"

# Row 9 - Confirm New Password
$ws1.Range("C9").Value = "Neues Passwort bestätigen"
$ws1.Range("E9").Value = "Confirmar nueva contraseña"

# Row 10 - Change Password
$ws1.Range("C10").Value = "Passwort ändern"
$ws1.Range("E10").Value = "Cambiar contraseña"

# Row 11 - Password recover (D gets an odd placeholder per source data)
$ws1.Range("D11").Value = "Readeooooo messageiooooo"

# Row 13 - Read message
$ws1.Range("C13").Value = "Nachricht lesen"
$ws1.Range("E13").Value = "Leer mensaje"

# Row 14 - Are you sure you want to delete this message?
$ws1.Range("C14").Value = "Sind Sie sicher, dass Sie diese Nachricht löschen möchten?"
$ws1.Range("E14").Value = "¿Está seguro que desea eliminar este mensaje?"

# Row 15 - Search in messages
$ws1.Range("C15").Value = "In Nachrichten suchen"
$ws1.Range("E15").Value = "Buscar en mensajes"

# Row 16 - Compose Message
$ws1.Range("C16").Value = "Nachricht verfassen"
$ws1.Range("E16").Value = "Redactar mensaje"

# Row 17 - No date provided
$ws1.Range("C17").Value = "Kein Datum angegeben"
$ws1.Range("E17").Value = "No hay fecha proporcionadada"

# Row 18 - Quick Start
$ws1.Range("C18").Value = "Schnellstart"
$ws1.Range("E18").Value = "Inicio rápido"

# ---------------------------------------------------------------------
# Sheet "Sheet1": same column layout (A=Identifier label, B=English,
# C=German, D=Portuguese, E=Spanish, F=French). Rows 2-6 get new
# translations (row 1 header text is unchanged).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")

# Row 2 - Text for translation 1
$ws2.Range("C2").Value = " Traduction 1"
$ws2.Range("E2").Value = "Traducción 1"
$ws2.Range("F2").Value = " Traducción 1"

# Row 3 - Text for translation 2
$ws2.Range("C3").Value = "Traduction 2"
$ws2.Range("E3").Value = " Traducción 2"

# Row 4 - Text for translation 3
$ws2.Range("C4").Value = "Traduction 3 "
$ws2.Range("E4").Value = "Traducción 3"

# Row 5 - Text for translation 4
$ws2.Range("C5").Value = "Traduction 4"
$ws2.Range("E5").Value = "Traducción 4"

# Row 6 - Text for translation 5
$ws2.Range("C6").Value = "Traduction 5"
$ws2.Range("E6").Value = "Traducción 5"
